$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Q3 ranking table (rows 24-28), values entered in the order that
# produces the same shared-string ordering as the authored workbook:
# "Lower price point", then the question text, then the other two choices.
$ws.Range("A26").Value = "Lower price point"
$ws.Range("A24").Value = "Q3. How would you improve your current product (Rank 1 to 3)"
$ws.Range("A27").Value = "Better package"
$ws.Range("A28").Value = "Improved cleaning"

$ws.Range("A25").Value = "Answer Choices"
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3

$ws.Range("B26").Value = 91
$ws.Range("C26").Value = 33
$ws.Range("D26").Value = 19

$ws.Range("B27").Value = 31
$ws.Range("C27").Value = 37
$ws.Range("D27").Value = 16

$ws.Range("B28").Value = 18
$ws.Range("C28").Value = 11
$ws.Range("D28").Value = 17

# Match the saved selection state: whole row 21 selected, active cell A21.
$ws.Application.GoTo($ws.Range("A21:XFD21"))
